$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.01385889889611
$ws.Range("D2").Value = 1.016135689362683
$ws.Range("E2").Value = 1.015692096259334
$ws.Range("F2").Value = 1.025281038456786
$ws.Range("I2").Value = 1.025701409304815
$ws.Range("J2").Value = 1.019093454337336
$ws.Range("K2").Value = 1.018988763173778
$ws.Range("L2").Value = 1.018546495265856
$ws.Range("M2").Value = 1.028107069382837
$ws.Range("N2").Value = 1.020540683934148
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.015116750629211
$ws.Range("D3").Value = 1.017239546257658
$ws.Range("E3").Value = 1.016768840607863
$ws.Range("F3").Value = 1.026649954282592
$ws.Range("I3").Value = 1.025687028427521
$ws.Range("J3").Value = 1.019983996413019
$ws.Range("K3").Value = 1.019897107584716
$ws.Range("L3").Value = 1.019427705287252
$ws.Range("M3").Value = 1.029281730353043
$ws.Range("N3").Value = 1.021432490681726
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.015930543923489
$ws.Range("D4").Value = 1.017953988048463
$ws.Range("E4").Value = 1.017465767397186
$ws.Range("F4").Value = 1.027533958980057
$ws.Range("I4").Value = 1.025675215049564
$ws.Range("J4").Value = 1.020559692438234
$ws.Range("K4").Value = 1.020484462555606
$ws.Range("L4").Value = 1.019997524941466
$ws.Range("M4").Value = 1.030039523664973
$ws.Range("N4").Value = 1.022009004261329
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.016272636612942
$ws.Range("D5").Value = 1.018254382472046
$ws.Range("E5").Value = 1.017758805223073
$ws.Range("F5").Value = 1.027905172278182
$ws.Range("I5").Value = 1.025669647510118
$ws.Range("J5").Value = 1.020801586448179
$ws.Range("K5").Value = 1.020731291359345
$ws.Range("L5").Value = 1.020236987164904
$ws.Range("M5").Value = 1.030357553620779
$ws.Range("N5").Value = 1.022251241788514
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.016330074040484
$ws.Range("D6").Value = 1.018304822596316
$ws.Range("E6").Value = 1.017808010531424
$ws.Range("F6").Value = 1.027967475941239
$ws.Range("I6").Value = 1.02566867741845
$ws.Range("J6").Value = 1.02084219402484
$ws.Range("K6").Value = 1.020772729464924
$ws.Range("L6").Value = 1.020277188694974
$ws.Range("M6").Value = 1.030410920271507
$ws.Range("N6").Value = 1.022291907032591
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.015935115081327
$ws.Range("D7").Value = 1.017958001761928
$ws.Range("E7").Value = 1.01746968278525
$ws.Range("F7").Value = 1.027538920806318
$ws.Range("I7").Value = 1.025675143019131
$ws.Range("J7").Value = 1.020562925141062
$ws.Range("K7").Value = 1.020487761066327
$ws.Range("L7").Value = 1.020000724999018
$ws.Range("M7").Value = 1.030043775341176
$ws.Range("N7").Value = 1.022012241554966
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.014284022573849
$ws.Range("D8").Value = 1.016508708312211
$ws.Range("E8").Value = 1.01605594638177
$ws.Range("F8").Value = 1.025744039238568
$ws.Range("I8").Value = 1.025697068052583
$ws.Range("J8").Value = 1.019394530913413
$ws.Range("K8").Value = 1.019295827183099
$ws.Range("L8").Value = 1.018844384368643
$ws.Range("M8").Value = 1.02850452716622
$ws.Range("N8").Value = 1.020842188073498
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.011373529327463
$ws.Range("D9").Value = 1.013956109407429
$ws.Range("E9").Value = 1.013566229274661
$ws.Range("F9").Value = 1.022567506344935
$ws.Range("I9").Value = 1.025716528129859
$ws.Range("J9").Value = 1.017331426342908
$ws.Range("K9").Value = 1.017192314837466
$ws.Range("L9").Value = 1.016803768605725
$ws.Range("M9").Value = 1.025774534417281
$ws.Range("N9").Value = 1.01877615365786
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.00943230392377
$ws.Range("D10").Value = 1.012255096638489
$ws.Range("E10").Value = 1.011907298510126
$ws.Range("F10").Value = 1.020440412105491
$ws.Range("I10").Value = 1.02571665357834
$ws.Range("J10").Value = 1.01595306051072
$ws.Range("K10").Value = 1.015787734903381
$ws.Range("L10").Value = 1.015441245273222
$ws.Range("M10").Value = 1.023942540210484
$ws.Range("N10").Value = 1.01739583038809
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.008591477314588
$ws.Range("D11").Value = 1.011518683774729
$ws.Range("E11").Value = 1.01118914776302
$ws.Range("F11").Value = 1.019517084480588
$ws.Range("I11").Value = 1.025713669424727
$ws.Range("J11").Value = 1.015355485260069
$ws.Range("K11").Value = 1.015178981266717
$ws.Range("L11").Value = 1.014850733889139
$ws.Range("M11").Value = 1.023146388155564
$ws.Range("N11").Value = 1.01679740651204
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.008279113638703
$ws.Range("D12").Value = 1.01124516546235
$ws.Range("E12").Value = 1.010922419046019
$ws.Range("F12").Value = 1.01917377293353
$ws.Range("I12").Value = 1.025712105138288
$ws.Range("J12").Value = 1.015133406739065
$ws.Range("K12").Value = 1.014952776578803
$ws.Range("L12").Value = 1.014631310044163
$ws.Range("M12").Value = 1.022850225206983
$ws.Range("N12").Value = 1.016575012614063
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.008346118723462
$ws.Range("D13").Value = 1.01130383528252
$ws.Range("E13").Value = 1.010979632190136
$ws.Range("F13").Value = 1.019247430168653
$ws.Range("I13").Value = 1.025712461300313
$ws.Range("J13").Value = 1.015181048460908
$ws.Range("K13").Value = 1.015001302200554
$ws.Range("L13").Value = 1.014678380930146
$ws.Range("M13").Value = 1.022913772971066
$ws.Range("N13").Value = 1.016622721992616
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.008565658139098
$ws.Range("D14").Value = 1.011496074307798
$ws.Range("E14").Value = 1.011167099398026
$ws.Range("F14").Value = 1.01948871333245
$ws.Range("I14").Value = 1.025713549415551
$ws.Range("J14").Value = 1.01533713048745
$ws.Range("K14").Value = 1.015160284893406
$ws.Range("L14").Value = 1.014832597922559
$ws.Range("M14").Value = 1.023121916158196
$ws.Range("N14").Value = 1.016779025673539
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.008700917778965
$ws.Range("D15").Value = 1.011614521458901
$ws.Range("E15").Value = 1.011282607305943
$ws.Range("F15").Value = 1.019637329818308
$ws.Range("I15").Value = 1.025714159455582
$ws.Range("J15").Value = 1.015433282805669
$ws.Range("K15").Value = 1.015258227852966
$ws.Range("L15").Value = 1.014927605224451
$ws.Range("M15").Value = 1.023250102098721
$ws.Range("N15").Value = 1.016875314539079
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.009488100956452
$ws.Range("D16").Value = 1.012303972618587
$ws.Range("E16").Value = 1.011954963358158
$ws.Range("F16").Value = 1.020501641925814
$ws.Range("I16").Value = 1.025716787666102
$ws.Range("J16").Value = 1.015992703962863
$ws.Range("K16").Value = 1.015828123868239
$ws.Range("L16").Value = 1.015480424233957
$ws.Range("M16").Value = 1.023995317133214
$ws.Range("N16").Value = 1.017435530138483
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.009981807097589
$ws.Range("D17").Value = 1.012736482108398
$ws.Range("E17").Value = 1.012376760458358
$ws.Range("F17").Value = 1.021043188546163
$ws.Range("I17").Value = 1.02571762325623
$ws.Range("J17").Value = 1.016343416059765
$ws.Range("K17").Value = 1.016185452796266
$ws.Range("L17").Value = 1.015827049723991
$ws.Range("M17").Value = 1.024461995961223
$ws.Range("N17").Value = 1.01778674028679
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.010269752393111
$ws.Range("D18").Value = 1.012988770944871
$ws.Range("E18").Value = 1.012622804818871
$ws.Range("F18").Value = 1.021358843625142
$ws.Range("I18").Value = 1.02571781743024
$ws.Range("J18").Value = 1.016547909848174
$ws.Range("K18").Value = 1.016393822708802
$ws.Range("L18").Value = 1.016029179529442
$ws.Range("M18").Value = 1.024733923526996
$ws.Range("N18").Value = 1.017991524479833
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.010367930143972
$ws.Range("D19").Value = 1.013074797256448
$ws.Range("E19").Value = 1.012706702578356
$ws.Range("F19").Value = 1.02146643666258
$ws.Range("I19").Value = 1.025717833884442
$ws.Range("J19").Value = 1.016617625029817
$ws.Range("K19").Value = 1.016464862368713
$ws.Range("L19").Value = 1.016098091983625
$ws.Range("M19").Value = 1.024826596683905
$ws.Range("N19").Value = 1.018061338665031
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.009928839729809
$ws.Range("D20").Value = 1.012690076568687
$ws.Range("E20").Value = 1.012331503854091
$ws.Range("F20").Value = 1.020985108493243
$ws.Range("I20").Value = 1.025717563928957
$ws.Range("J20").Value = 1.016305795295009
$ws.Range("K20").Value = 1.016147120386655
$ws.Range("L20").Value = 1.015789865377881
$ws.Range("M20").Value = 1.024411954548323
$ws.Range("N20").Value = 1.017749066096232
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.008501010476306
$ws.Range("D21").Value = 1.011439464222446
$ws.Range("E21").Value = 1.011111894320532
$ws.Range("F21").Value = 1.019417671051431
$ws.Range("I21").Value = 1.025713241570915
$ws.Range("J21").Value = 1.015291171343441
$ws.Range("K21").Value = 1.015113470857648
$ws.Range("L21").Value = 1.01478718713164
$ws.Range("M21").Value = 1.023060635254451
$ws.Range("N21").Value = 1.016733001262273
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.007603025568567
$ws.Range("D22").Value = 1.010653258809998
$ws.Range("E22").Value = 1.010345217750034
$ws.Range("F22").Value = 1.018430153972539
$ws.Range("I22").Value = 1.025707886903917
$ws.Range("J22").Value = 1.014652585189599
$ws.Range("K22").Value = 1.014463072933887
$ws.Range("L22").Value = 1.014156290255036
$ws.Range("M22").Value = 1.022208479498067
$ws.Range("N22").Value = 1.016093508242847
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.008079089481204
$ws.Range("D23").Value = 1.011070032086333
$ws.Range("E23").Value = 1.010751634933945
$ws.Range("F23").Value = 1.018953846956906
$ws.Range("I23").Value = 1.025710975258589
$ws.Range("J23").Value = 1.014991174323436
$ws.Range("K23").Value = 1.014807909434654
$ws.Range("L23").Value = 1.014490786167564
$ws.Range("M23").Value = 1.022660463983661
$ws.Range("N23").Value = 1.01643257821209
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.009952773480905
$ws.Range("D24").Value = 1.012711045193251
$ws.Range("E24").Value = 1.012351953312211
$ws.Range("F24").Value = 1.021011353048413
$ws.Range("I24").Value = 1.025717591642484
$ws.Range("J24").Value = 1.016322794717463
$ws.Range("K24").Value = 1.016164441319388
$ws.Range("L24").Value = 1.015806667540792
$ws.Range("M24").Value = 1.024434566968578
$ws.Range("N24").Value = 1.017766089659816
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.012126105971104
$ws.Range("D25").Value = 1.014615882081278
$ws.Range("E25").Value = 1.014209717816454
$ws.Range("F25").Value = 1.023390360755825
$ws.Range("I25").Value = 1.025713764933832
$ws.Range("J25").Value = 1.017865302515839
$ws.Range("K25").Value = 1.017736509969724
$ws.Range("L25").Value = 1.017331681829682
$ws.Range("M25").Value = 1.026482407626778
$ws.Range("N25").Value = 1.019310787996194
